# Slide 1, shape "Textplatzhalter 3" (the byline paragraph), reorder the
# author names from "Jérôme Imfeld, Julian Schuhmacher" to
# "Julian Schuhmacher, Jérôme Imfeld" while keeping "by ... and Joel Fimian".
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Original run layout inside paragraph 4 ("by Jérôme Imfeld, Julian Schuhmacher and Joel Fimian"):
#   "by" | " Jérôme " | "Imfeld" | ", Julian Schuhmacher " | "and" | " Joel Fimian"

# 1) Replace the " Jérôme " run's text in place -> " Julian Schuhmacher, "
#    (keeps that run's original formatting/rPr untouched)
$runJerome = $tr.Characters(38, 8)
$runJerome.Text = " Julian Schuhmacher, "

# 2) Insert a brand-new run right after it containing "Jérôme " with the
#    same character formatting as its neighbour, then force an explicit
#    font property so the engine keeps it as its own distinct run/rPr
#    instead of silently merging it back into the previous run.
$afterRun = $tr.Characters(38, 21)
$newRun = $afterRun.InsertAfter([string][char]0x4A + [char]0xE9 + "r" + [char]0xF4 + "me ")
$newRun.Font.Size = 18

# 3) The old ", Julian Schuhmacher " run (now shifted right because of the
#    edits above) collapses down to a single space.
$runOldJulian = $tr.Characters(72, 21)
$runOldJulian.Text = " "
